{"js": "const replacements = [\n  [\"32\u00f76=\", \"52\u00f73=\"],\n  [\"99\u00f72=\", \"80\u00f72=\"],\n  [\"88\u00f76=\", \"78\u00f77=\"],\n  [\"84\u00f77=\", \"74\u00f72=\"],\n  [\"35\u00f76=\", \"54\u00f76=\"],\n  [\"69\u00f77=\", \"72\u00f79=\"],\n  [\"77\u00f74=\", \"75\u00f74=\"],\n  [\"65\u00f74=\", \"64\u00f78=\"],\n  [\"33\u00f78=\", \"82\u00f76=\"],\n  [\"91\u00f77=\", \"36\u00f74=\"],\n  [\"32\u00f73=\", \"17\u00f75=\"],\n  [\"78\u00f75=\", \"44\u00f77=\"],\n  [\"40\u00f76=\", \"48\u00f79=\"],\n  [\"51\u00f78=\", \"28\u00f77=\"],\n  [\"27\u00f73=\", \"99\u00f77=\"],\n  [\"71\u00f74=\", \"85\u00f74=\"],\n  [\"44\u00f76=\", \"20\u00f75=\"],\n  [\"24\u00f73=\", \"11\u00f77=\"],\n  [\"75\u00f76=\", \"74\u00f77=\"],\n  [\"29\u00f75=\", \"54\u00f78=\"],\n  [\"89\u00f74=\", \"62\u00f79=\"],\n  [\"45\u00f72=\", \"15\u00f76=\"],\n  [\"25\u00f72=\", \"58\u00f78=\"],\n  [\"60\u00f73=\", \"32\u00f79=\"],\n  [\"12\u00f74=\", \"26\u00f73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"32\u00f76=\", \"52\u00f73=\"),\n    @(\"99\u00f72=\", \"80\u00f72=\"),\n    @(\"88\u00f76=\", \"78\u00f77=\"),\n    @(\"84\u00f77=\", \"74\u00f72=\"),\n    @(\"35\u00f76=\", \"54\u00f76=\"),\n    @(\"69\u00f77=\", \"72\u00f79=\"),\n    @(\"77\u00f74=\", \"75\u00f74=\"),\n    @(\"65\u00f74=\", \"64\u00f78=\"),\n    @(\"33\u00f78=\", \"82\u00f76=\"),\n    @(\"91\u00f77=\", \"36\u00f74=\"),\n    @(\"32\u00f73=\", \"17\u00f75=\"),\n    @(\"78\u00f75=\", \"44\u00f77=\"),\n    @(\"40\u00f76=\", \"48\u00f79=\"),\n    @(\"51\u00f78=\", \"28\u00f77=\"),\n    @(\"27\u00f73=\", \"99\u00f77=\"),\n    @(\"71\u00f74=\", \"85\u00f74=\"),\n    @(\"44\u00f76=\", \"20\u00f75=\"),\n    @(\"24\u00f73=\", \"11\u00f77=\"),\n    @(\"75\u00f76=\", \"74\u00f77=\"),\n    @(\"29\u00f75=\", \"54\u00f78=\"),\n    @(\"89\u00f74=\", \"62\u00f79=\"),\n    @(\"45\u00f72=\", \"15\u00f76=\"),\n    @(\"25\u00f72=\", \"58\u00f78=\"),\n    @(\"60\u00f73=\", \"32\u00f79=\"),\n    @(\"12\u00f74=\", \"26\u00f73=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
